# Remove the "A、" / "B、" / "C、" / "D、" option-letter prefixes from the
# four multiple-choice answer sets that still carried them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D92").Value = "心室纤维性颤动"
$ws.Range("E92").Value = "心房纤维性颤动"
$ws.Range("F92").Value = "心博停止"
$ws.Range("G92").Value = "以上都可以"

$ws.Range("D96").Value = "紧缩压榨感"
$ws.Range("E96").Value = "钝痛"
$ws.Range("F96").Value = "锐痛"
$ws.Range("G96").Value = "持续痛"

$ws.Range("D99").Value = "前脱位"
$ws.Range("E99").Value = "后脱位"
$ws.Range("F99").Value = "下脱位"
$ws.Range("G99").Value = "盂上脱位"

$ws.Range("D101").Value = "冲--脱--泡--盖--送"
$ws.Range("E101").Value = "泡--盖--送--冲--脱"
$ws.Range("F101").Value = "冲--泡--脱--盖--送"
$ws.Range("G101").Value = "泡--脱--冲--盖--送"

# The shortened text in row 92 no longer needs the wrapped two-line row
# height, so re-autofit that row (drops the explicit ht="24").
$ws.Rows.Item(92).AutoFit()

# Selection moved from C5 (with the view scrolled so A19 was the
# top-left visible cell) to E6, with the view back at the natural top.
$ws.Range("E6").Select()
